$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (A8): change from the "last row" date format (YYYY-MM-DD) to the
# regular date-time format (YYYY-MM-DD HH:MM:SS) used by the rows above it.
$ws.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 9: becomes the new "last row" with the special date-only format.
$ws.Range("A9").Value = 44515
$ws.Range("A9").NumberFormat = "YYYY-MM-DD"
$ws.Range("B9").Value = 57871.5
